$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume cells keep their original text formatting so values
# such as "1.001" or "30.843.44" are stored as literal text, not coerced to numbers/dates.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "30.843.44"
$ws.Range("E2").Value = "  +2.04%  "

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.04"
$ws.Range("E3").Value = "  +1.33%  "

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "245.02"
$ws.Range("E5").Value = "  +4.38%  "

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4786"
$ws.Range("E7").Value = "  +1.79%  "

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2911"
$ws.Range("E8").Value = "  +2.27%  "

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "42.89"
$ws.Range("E9").Value = "  +2.81%  "

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06581"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "21.41"
$ws.Range("E11").Value = "  +1.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07791"

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.906.47"
$ws.Range("E13").Value = "  +2.28%  "

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "97.26"
$ws.Range("E14").Value = "  +0.76%  "

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7406"
$ws.Range("E15").Value = "  +7.24%  "

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "5.195"
$ws.Range("E16").Value = "  +2.13%  "

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "281.80"
$ws.Range("E17").Value = "  +5.78%  "

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "30.835.99"
$ws.Range("E18").Value = "  +2.04%  "

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "13.62"
$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007643"
$ws.Range("E20").Value = "  -0.59%  "

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9993"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "2.163.35"
$ws.Range("E22").Value = "  +2.76%  "

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "5.323"
$ws.Range("E23").Value = "  +1.70%  "

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "6.258"
$ws.Range("E25").Value = "  +1.80%  "

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "9.394"
$ws.Range("E26").Value = "  -0.93%  "

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "166.46"
$ws.Range("E27").Value = "  +0.47%  "

$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = "19.19"
$ws.Range("E28").Value = "  +2.35%  "

$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D29").Value = "1.994"
$ws.Range("E29").Value = "  +3.26%  "

$ws.Range("D30:E30").NumberFormat = "@"
$ws.Range("D30").Value = "1.383"
$ws.Range("E30").Value = "  +0.91%  "

$ws.Range("D31:E31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1004"
$ws.Range("E31").Value = "  +1.32%  "

$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("D32").Value = "1.521"
$ws.Range("E32").Value = "  +4.60%  "

$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("D33").Value = "4.399"
$ws.Range("E33").Value = "  +1.24%  "

$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("D34").Value = "4.144"
$ws.Range("E34").Value = "  +2.56%  "

$ws.Range("D35:E35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04792"
$ws.Range("E35").Value = "  +1.14%  "

$ws.Range("D36:E36").NumberFormat = "@"
$ws.Range("D36").Value = "1.134"
$ws.Range("E36").Value = "  +0.44%  "

$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7085"
$ws.Range("E37").Value = "  +1.40%  "

$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "2.720"
$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01877"
$ws.Range("E39").Value = "  +0.71%  "

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "2.775"
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "6.460"
$ws.Range("E41").Value = "  +3.28%  "

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "71.07"
$ws.Range("E42").Value = "  -2.00%  "

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "1.938"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4230"
$ws.Range("E44").Value = "  +1.84%  "

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8492"
$ws.Range("E45").Value = "  +2.07%  "

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "102.79"
$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "9.528"
$ws.Range("E48").Value = "  +4.33%  "

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "7.179"
$ws.Range("E49").Value = "  +1.44%  "

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "950.13"
$ws.Range("E50").Value = "  -2.86%  "

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "35.40"
$ws.Range("E51").Value = "  +2.72%  "
